# Update the opening paragraph:
#   "This is a Microsoft word document." -> adds two trailing spaces, then
#   appends a red (FF0000) parenthetical note, typed as three separate runs:
#     "(This is a change – Ve" / "rsion for main branch" / ")"

$d = $word.ActiveDocument

# Add the two trailing spaces to the existing sentence (keeps it in its own run).
$d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false, $false,
    $true, 1, $false, "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs(1).Range

# Insert the new text in three separate chunks (mirrors how the edit was
# authored), coloring each chunk's own range right after insertion so the
# runs stay distinct instead of being coalesced into one big red run.

$pos0 = $p1.End - 1
$p1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$pos1 = $p1.End - 1
$d.Range($pos0, $pos1).Font.Color = 255

$p1.InsertAfter("rsion for main branch")
$pos2 = $p1.End - 1
$d.Range($pos1, $pos2).Font.Color = 255

$p1.InsertAfter(")")
$pos3 = $p1.End - 1
$d.Range($pos2, $pos3).Font.Color = 255
